$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-CAN_CurrentLo")

# Note: a leading apostrophe is used when assigning text values so that the
# cells' existing "quote prefix" text formatting (and thus their original
# cell style) is preserved instead of being dropped by a plain text write.

# Row 44: Amplifier gain resistors changed from 62 kOhms to 56 kOhms
$ws.Range("B44").Value = "'CRCW060356K0FKEA"
$ws.Range("E44").Value = "'56 kOhms"
$ws.Range("F44").Value = "'Thin Film Resistor 56 kOhms 1%"

# Row 46: R26 removed from the POT BOURNS TRIMPOT 3224W group, quantity 6 -> 5
$ws.Range("A46").Value = "'R13, R14, R18, R24, R25"
$ws.Range("D46").Value = 5

# Row 51: R26 added to the 4.7 kOhms resistor group, quantity 3 -> 4
$ws.Range("A51").Value = "'R21, R22, R23, R26"
$ws.Range("D51").Value = 4
